# Update view-count figures in column F (reflecting a fresh scrape) across
# the four sheets: 展览, 演出, 本地生活, 全部类型.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        2 = 2697; 3 = 586; 4 = 473; 8 = 1256; 9 = 589; 10 = 321; 11 = 12;
        12 = 140; 13 = 380; 14 = 5882; 15 = 102; 16 = 1824; 17 = 4314;
        18 = 447; 21 = 5027; 22 = 6481; 25 = 709; 26 = 3846; 27 = 514;
        31 = 1007; 32 = 1437; 33 = 503; 34 = 606; 35 = 1627; 36 = 214;
        37 = 1771; 40 = 1340; 41 = 645; 42 = 103; 43 = 3507; 46 = 419;
        47 = 11; 48 = 52
    }
    "演出" = @{
        10 = 17; 11 = 17
    }
    "本地生活" = @{
        2 = 4082
    }
    "全部类型" = @{
        2 = 4082; 3 = 2697; 4 = 586; 5 = 473; 12 = 1256; 13 = 17; 14 = 589;
        15 = 321; 16 = 140; 17 = 380; 18 = 102; 19 = 1825; 20 = 4314;
        21 = 5027; 22 = 5027; 25 = 709; 26 = 3846; 27 = 515; 30 = 1007;
        31 = 1437; 32 = 503; 33 = 606; 34 = 1627; 35 = 214; 36 = 1771;
        40 = 645; 42 = 103; 44 = 3507; 48 = 52
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
